$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.811.86"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "3.576.76"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "578.53"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.41%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "188.97"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.632"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("D8").Value = "3.573.30"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -0.59%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "55.71"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("E13").Value = "  +1.98%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "9.62"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "4.155.12"
$ws.Range("E15").Value = "  -1.20%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "19.86"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "3.579.46"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "69.778.32"
$ws.Range("E18").Value = "  -1.23%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.64"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -1.16%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "475.04"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -4.26%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "19.46"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +12.29%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "5.04"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -6.98%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "96.04"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +5.35%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "4.38"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.84%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "3.00"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -4.03%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "11.02"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.13%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.37"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.26%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.90"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +4.39%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "32.39"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +0.51%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "66.11"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "582.23"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -6.34%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "39.02"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "0.0₃0797"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("E39").Value = "  -4.12%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.24"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +18.80%  "
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("E42").Value = "  -6.20%  "
$ws.Range("D43").Value = "3.237.99"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("E44").Value = "  -6.19%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").Value = "  -0.37%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.41"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  -0.28%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -5.63%  "
